$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.503.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.133.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.85%  "

$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5253"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4561"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.45"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09120"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.189"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.137.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.885"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001168"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.42%  "

$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06720"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.31%  "

$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.373"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.596.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.380"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.389.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.597"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.222"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.715"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.36%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1082"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.388"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.033"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.142"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02640"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06984"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2354"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6977"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.273"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.16%  "

$ws.Range("E44").Value = "  +5.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6504"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.345"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000373"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.709"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.248"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07285"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.09%  "
